# Adding data summaries: append the 2020-04-05 row to the "Confirmados"
# and "Mortes" sheets (sheet1 / sheet2). The "Populacao" sheet is left
# untouched, its data does not change.

$wb = $excel.ActiveWorkbook

$wsConfirmados = $wb.Worksheets.Item("Confirmados")
$wsMortes = $wb.Worksheets.Item("Mortes")

$newRow = 42
$dateLabel = "2020-04-05"

# Values for the new row, indexed by column number (B=2 .. AB=28),
# in the same state order as every other row in the sheet.
$confirmados = @{
    2 = 48; 3 = 28; 4 = 29; 5 = 417; 6 = 401; 7 = 823; 8 = 468; 9 = 166;
    10 = 115; 11 = 96; 12 = 60; 13 = 65; 14 = 498; 15 = 86; 16 = 34;
    17 = 438; 18 = 201; 19 = 23; 20 = 1394; 21 = 242; 22 = 418; 23 = 12;
    24 = 42; 25 = 357; 26 = 4620; 27 = 32; 28 = 17
}

$mortes = @{
    2 = 0; 3 = 2; 4 = 2; 5 = 14; 6 = 9; 7 = 26; 8 = 7; 9 = 6;
    10 = 3; 11 = 2; 12 = 1; 13 = 1; 14 = 6; 15 = 1; 16 = 4;
    17 = 9; 18 = 21; 19 = 4; 20 = 64; 21 = 7; 22 = 7; 23 = 1;
    24 = 1; 25 = 10; 26 = 275; 27 = 3; 28 = 0
}

function Add-DateRow {
    param($ws, $rowValues)

    # Force column A to be stored as text so that the date-like string
    # "2020-04-05" is kept as a label (matching the rest of the column)
    # rather than being converted into a date serial number.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dateLabel
    $cellA.Style = "Normal"

    foreach ($col in $rowValues.Keys) {
        $ws.Cells.Item($newRow, [int]$col).Value = $rowValues[$col]
    }
}

Add-DateRow $wsConfirmados $confirmados
Add-DateRow $wsMortes $mortes

Write-Output "Added row $newRow ($dateLabel) to Confirmados and Mortes sheets"
